$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: D8 gets a plain-text note
$ws.Range("D8").Value = 'Co nam moduły dostarczą :)'

# Row 9: D9 gets a long rich-text explanation (left/center aligned, wrapped)
$ws.Rows("9").RowHeight = 195
$ws.Range("D9").Value = 'Prosimy o precyzowanie przez modyły jaki będzie format danych wyjściowych DOKŁADNIE!! np. przetworzony cały sygnał - Vector<double> (czyli amplitudy kolejnych próbek), parametry jakieśtam do wyświetlenia w tabeli - map<Qstring,double> (gdzie Qstring to nazwa parametru a double to jego wartość), jeśli to dane do histogramu - poprosimy Vector<double>, jako wartości kolejnych słupków i jakiś double, który będzie szerokością słupka,  jeżeli to np. wykryte załamki to poprosimy Vector<double> - czyli numery tychże próbek, dla których stwierdzono, że są załamkiem, jeżeli na wykresie ma być zaznaczony jakiś odcinek to najlepiej Vector<double>, gdzie kolejne wartości to będzie nr próbki startowej, długość przedziału, nr kolejnej próbki startowej, długość przedziału itd. Sorry, same nazwy tego, co będzie zwracane nic nie mówią....'
$ws.Range("D9").HorizontalAlignment = -4131
$ws.Range("D9").VerticalAlignment = -4108
$ws.Range("D9").WrapText = $true
$ws.Range("D9").Characters(1,7).Font.Bold = $true
$ws.Range("D9").Characters(118,14).Font.Bold = $true
$ws.Range("D9").Characters(215,20).Font.Bold = $true
$ws.Range("D9").Characters(337,14).Font.Bold = $true
$ws.Range("D9").Characters(393,6).Font.Bold = $true
$ws.Range("D9").Characters(478,14).Font.Bold = $true
$ws.Range("D9").Characters(629,14).Font.Bold = $true
$ws.Range("D9").Characters(8,110).Font.Bold = $false
$ws.Range("D9").Characters(132,83).Font.Bold = $false
$ws.Range("D9").Characters(235,102).Font.Bold = $false
$ws.Range("D9").Characters(351,42).Font.Bold = $false
$ws.Range("D9").Characters(399,79).Font.Bold = $false
$ws.Range("D9").Characters(492,137).Font.Bold = $false
$ws.Range("D9").Characters(643,190).Font.Bold = $false

# Row 32: taller row; E32 gets wrap text (D32 reuses the existing wrap style)
$ws.Rows("32").RowHeight = 30
$ws.Range("D32").WrapText = $true
$ws.Range("E32").WrapText = $true

# Row 9: B9/C9 get "vertical top" alignment
$ws.Range("B9").VerticalAlignment = -4160
$ws.Range("C9").VerticalAlignment = -4160

# Selection moves to A9 (matches the last-saved cursor position)
$ws.Range("A9").Select()

# Page setup: portrait printing, paper size 9 were saved with the sheet
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
